$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace worker names in column A (rows 2-12) with generic "worker N" labels
for ($i = 1; $i -le 11; $i++) {
    $row = $i + 1
    $ws.Range("A$row").Value = "worker $i"
}

# Update the active selection to match the author's final cursor position
$ws.Range("B18").Select()
